# Updated symbol list - refresh Price (D) and Volume(1h) (E) columns
# Values are stored as literal text (matching source data export format),
# so we use the leading apostrophe to force text entry, same as typing
# a number-looking value into a cell and prefixing it with ' in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.69"
$ws.Range("E2").Value = "'-0.80%"
$ws.Range("D3").Value = "'37.50"
$ws.Range("E3").Value = "'7.53%"
$ws.Range("D4").Value = "'5.009"
$ws.Range("E4").Value = "'-2.29%"
$ws.Range("D5").Value = "'0.07853"
$ws.Range("E5").Value = "'1.11%"
$ws.Range("D6").Value = "'2.265"
$ws.Range("E6").Value = "'-3.86%"
$ws.Range("D7").Value = "'8.027"
$ws.Range("E7").Value = "'0.09%"
$ws.Range("E8").Value = "'2.08%"
$ws.Range("E9").Value = "'-1.63%"
$ws.Range("D10").Value = "'0.1886"
$ws.Range("E10").Value = "'5.10%"
$ws.Range("D11").Value = "'0.09258"
$ws.Range("E11").Value = "'-7.87%"
$ws.Range("D12").Value = "'0.08478"
$ws.Range("E12").Value = "'-0.38%"
$ws.Range("D13").Value = "'0.03535"
$ws.Range("E13").Value = "'6.60%"
$ws.Range("D14").Value = "'0.09947"
$ws.Range("E14").Value = "'0.54%"
$ws.Range("D15").Value = "'0.001486"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("D16").Value = "'0.005630"
$ws.Range("E16").Value = "'-2.28%"
$ws.Range("E17").Value = "'0.03%"
$ws.Range("D18").Value = "'2.153"
$ws.Range("E18").Value = "'-0.48%"
$ws.Range("E19").Value = "'2.85%"
$ws.Range("E20").Value = "'-0.15%"
$ws.Range("D21").Value = "'4.770"
$ws.Range("E21").Value = "'11.14%"
$ws.Range("D22").Value = "'0.2201"
$ws.Range("E22").Value = "'-7.72%"
$ws.Range("D23").Value = "'0.04650"
$ws.Range("E23").Value = "'1.93%"
$ws.Range("D24").Value = "'0.001228"
$ws.Range("E24").Value = "'0.94%"
$ws.Range("D25").Value = "'0.004450"
$ws.Range("E25").Value = "'-0.41%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("D27").Value = "'0.0004744"
$ws.Range("E27").Value = "'28.30%"
$ws.Range("D39").Value = "'0.01767"
$ws.Range("E39").Value = "'-0.93%"
$ws.Range("D40").Value = "'0.04739"
$ws.Range("E40").Value = "'-0.27%"
$ws.Range("D41").Value = "'0.007901"
$ws.Range("E41").Value = "'1.95%"
$ws.Range("E42").Value = "'-1.55%"
$ws.Range("D43").Value = "'0.007652"
$ws.Range("E43").Value = "'8.20%"
$ws.Range("D44").Value = "'0.002228"
$ws.Range("E44").Value = "'5.94%"
$ws.Range("D45").Value = "'0.01018"
$ws.Range("E45").Value = "'6.98%"
$ws.Range("D46").Value = "'0.00006055"
$ws.Range("E46").Value = "'-1.00%"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("D48").Value = "'8.670"
$ws.Range("E48").Value = "'217.33%"
$ws.Range("D49").Value = "'0.002686"
$ws.Range("E49").Value = "'34.41%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.07%"
